# Generowany przebieg obiegow do excela
# Updates the "Termin" (column D) classification codes on sheet Arkusz1:
#   - trailing "-" suffix codes ("C-", "H-") are normalized to their base
#     letter ("C", "H")
#   - bracketed numeric codes that used to be stored as plain negative
#     numbers (-7, -1, -5) are now stored as quoted text labels ("[7]",
#     "[1]", "[5]")
#   - parenthesised range codes ("(1-6)", "(2-5)", "(1-4)") become
#     bracketed range codes ("[1-6]", "[2-5]", "[1-4]")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "C-" -> "C" ---------------------------------------------------
$cDashToC = @(8, 9, 63, 64, 101, 102, 115, 116, 121, 122, 123, 124, 125, 143, 144)
foreach ($r in $cDashToC) {
    $ws.Cells.Item($r, 4).Value = "C"
}

# --- "H-" -> "H" ---------------------------------------------------
$hDashToH = @(99, 100, 107, 108, 109, 110, 129, 132, 141, 148, 155)
foreach ($r in $hDashToH) {
    $ws.Cells.Item($r, 4).Value = "H"
}

# --- "B" stays "B" (shared-string table re-ordered upstream) -------
$ws.Cells.Item(49, 4).Value = "B"

# --- "(1-6)" -> "[1-6]" ----------------------------------------------
$parenOneSixToBracket = @(43, 44, 86, 94)
foreach ($r in $parenOneSixToBracket) {
    $ws.Cells.Item($r, 4).Value = "[1-6]"
}

# --- numeric codes (-7 / -1 / -5) become quoted bracket text, and the --
# --- remaining parenthesised range codes become bracketed range codes -
# A leading apostrophe forces Excel to store the value as text and marks
# the cell's style with quotePrefix="1", matching the target workbook.
$negSevenToBracket = @(82, 83, 87, 106)
foreach ($r in $negSevenToBracket) {
    $ws.Cells.Item($r, 4).Value = "'[7]"
}
$ws.Cells.Item(117, 4).Value = "'[1]"
$ws.Cells.Item(118, 4).Value = "[2-5]"
$ws.Cells.Item(119, 4).Value = "[1-4]"
$ws.Cells.Item(120, 4).Value = "'[5]"

# --- move the active selection to N17, matching the authored state --
$ws.Range("N17").Select()
